# Fix the algorithm/conditions on filtering the status of candidates.
# The "Active Candidates" report is re-derived with an updated filter,
# which changes several rows and adds 3 new candidate rows (the sheet
# grows from 24 data rows to 27 data rows: A1:F25 -> A1:F28).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target data for rows 2-28 (Job ID, Company, Job Title, Candidate, Status, Action Date serial)
$rows = @(
  @(580, "Legion Security", "Sales Engineer (US)", "Michael Saladino", "CV Sent", 45988),
  @(580, "Legion Security", "Sales Engineer (US)", "Michael Maizel", "2nd Interview", 45978),
  @(663, "Blink Ops", "Sales Engineer UK", "Kev Pyart", "3rd Interview", 45966),
  @(663, "Blink Ops", "Sales Engineer UK", "LLOYD WEBB", "1st Interview", 45994),
  @(663, "Blink Ops", "Sales Engineer UK", "James Todd", "CV Sent", 45981),
  @(663, "Blink Ops", "Sales Engineer UK", "Graham Rance", "2nd Interview", 45992),
  @(663, "Blink Ops", "Sales Engineer UK", "Alistair Macrae", "1st Interview", 45987),
  @(663, "Blink Ops", "Sales Engineer UK", "Adam Evans", "1st Interview", 45985),
  @(673, "Redwood Software", "SE UK", "Nicholas Vlatko", "4th Interview", 45987),
  @(707, "Dash0", "Sales Engineer EMEA (UK, Nordics, Benelux, Germany) x 2", "Harry Kimpel", "1st Interview", 45992),
  @(707, "Dash0", "Sales Engineer EMEA (UK, Nordics, Benelux, Germany) x 2", "Patrick Schrimpf", "1st Interview", 45993),
  @(730, "PointFive", "PointFive SE EST", "Matthew Hughes", "1st Interview", 45994),
  @(745, "Blink Ops", "Senior Sales Engineer (Mid-ATL)", "James Gaidos", "1st Interview", 45993),
  @(745, "Blink Ops", "Senior Sales Engineer (Mid-ATL)", "Jakub Nogalski", "1st Interview", 45982),
  @(768, "Adaptive6", "Senior Sales Engineer (US)", "Matthew Hughes", "CV Sent", 45986),
  @(768, "Adaptive6", "Senior Sales Engineer (US)", "Artur Sirota", "CV Sent", 45988),
  @(777, "Adaptive6", "SE Director", "Hozefa Bata", "CV Sent", 45974),
  @(788, "Mabl", "Partner Manager", "Alix Moreira", "2nd Interview", 45995),
  @(788, "Mabl", "Partner Manager", "Juan Echeverri", "1st Interview", 45992),
  @(790, "Allium", "Growth Marketing", "Cameron Bernard", "CV Sent", 45982),
  @(790, "Allium", "Growth Marketing", "SHAIL SHAH", "CV Sent", 45978),
  @(790, "Allium", "Growth Marketing", "David Lambert", "CV Sent", 45980),
  @(790, "Allium", "Growth Marketing", "Julia Nelson", "CV Sent", 45978),
  @(790, "Allium", "Growth Marketing", "KC Patrick", "CV Sent", 45980),
  @(791, "Adaptive6", "Head of Sales (US)", "Joseph Crowley", "1st Interview", 45987),
  @(799, "Legion Security", "CS1 Legion - Sales Engineer", "Michael Maizel", "CV Sent", 45966),
  @(826, "Legit Security", "VP of Sales", "Joseph Crowley", "1st Interview", 45992)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Length; $i++) {
    $row = $rows[$i]
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 6).NumberFormat = "YYYY-MM-DD HH:MM:SS"
}
